$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cells whose new text would otherwise be auto-converted to a number by Excel ---
# Force these ranges to Text format first so the values stay as strings, matching the source data type (inlineStr).
$textForceCells = @("D5", "D8", "D10", "D15", "D17", "D20", "D22", "D25", "D28", "D29", "D38", "D43", "D45", "D46", "D48", "D50")
foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D5').Value = '212.42'
$ws.Range('D8').Value = '0.249'
$ws.Range('D10').Value = '18.85'
$ws.Range('D15').Value = '0.521'
$ws.Range('D17').Value = '62.59'
$ws.Range('D20').Value = '202.78'
$ws.Range('D22').Value = '9.32'
$ws.Range('D25').Value = '144.81'
$ws.Range('D28').Value = '15.20'
$ws.Range('D29').Value = '6.60'
$ws.Range('D38').Value = '0.803'
$ws.Range('D43').Value = '0.782'
$ws.Range('D45').Value = '92.23'
$ws.Range('D46').Value = '1.52'
$ws.Range('D48').Value = '53.99'
$ws.Range('D50').Value = '0.409'

# --- Cells whose new text is safely non-numeric (keeps its string type automatically) ---
$ws.Range('D2').Value = '26.375.44'
$ws.Range('E2').Value = '  +0.76%  '
$ws.Range('D3').Value = '1.624.40'
$ws.Range('E3').Value = '  +1.40%  '
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('E5').Value = '  +0.20%  '
$ws.Range('E6').Value = '  -0.11%  '
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('E8').Value = '  +0.15%  '
$ws.Range('E9').Value = '  +0.35%  '
$ws.Range('E10').Value = '  +3.77%  '
$ws.Range('E11').Value = '  +0.63%  '
$ws.Range('D12').Value = '1.849.78'
$ws.Range('E12').Value = '  +1.40%  '
$ws.Range('D13').Value = '1.613.37'
$ws.Range('E13').Value = '  +0.78%  '
$ws.Range('E14').Value = '  +0.44%  '
$ws.Range('E15').Value = '  +0.80%  '
$ws.Range('D16').Value = '26.387.58'
$ws.Range('E16').Value = '  +0.81%  '
$ws.Range('E17').Value = '  +2.51%  '
$ws.Range('D18').Value = '0.0₃0727'
$ws.Range('E18').Value = '  +0.00%  '
$ws.Range('E19').Value = '  -0.15%  '
$ws.Range('E20').Value = '  -0.27%  '
$ws.Range('E22').Value = '  +0.69%  '
$ws.Range('E23').Value = '  +0.49%  '
$ws.Range('E24').Value = '  -3.56%  '
$ws.Range('E25').Value = '  +0.11%  '
$ws.Range('E26').Value = '  -0.08%  '
$ws.Range('E27').Value = '  -2.47%  '
$ws.Range('B28').Value = 'EthereumClassic'
$ws.Range('C28').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('E28').Value = '  -0.02%  '
$ws.Range('B29').Value = 'Cosmos'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('E29').Value = '  +1.06%  '
$ws.Range('E30').Value = '  +5.14%  '
$ws.Range('E31').Value = '  +0.59%  '
$ws.Range('E32').Value = '  +1.51%  '
$ws.Range('E33').Value = '  -0.13%  '
$ws.Range('E34').Value = '  +0.64%  '
$ws.Range('E35').Value = '  +2.28%  '
$ws.Range('D36').Value = '1.160.61'
$ws.Range('E36').Value = '  +1.71%  '
$ws.Range('E37').Value = '  +0.36%  '
$ws.Range('E38').Value = '  +2.32%  '
$ws.Range('E40').Value = '  -0.12%  '
$ws.Range('E41').Value = '  -0.03%  '
$ws.Range('E42').Value = '  +3.63%  '
$ws.Range('E43').Value = '  -0.07%  '
$ws.Range('D44').Value = '1.762.65'
$ws.Range('E44').Value = '  +1.40%  '
$ws.Range('E45').Value = '  +0.02%  '
$ws.Range('B46').Value = 'RenderToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('E46').Value = '  +1.40%  '
$ws.Range('B47').Value = 'BabyDogeCoin'
$ws.Range('C47').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D47').Value = '0.0₆0104'
$ws.Range('E47').Value = '  +9.51%  '
$ws.Range('E48').Value = '  -0.23%  '
$ws.Range('E49').Value = '  +0.21%  '
$ws.Range('E51').Value = '  -0.21%  '
